$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@("2025년도 한국통신학회 아이디어 경진대회(25.02.06.~07.) 참가 지원", "관리자", "2025-01-09", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=29")
    ,@("[융합트랙] 2025학년도 스마트헬스케어SW 융합트랙 이수 신청 안내(학과: 보건행정학과)", "관리자", "2025-01-08", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=28")
    ,@("[융합트랙] 2025학년도 AI기반 디지털정밀농업 융합트랙 이수 신청 안내(학과: 식물자원학과)", "관리자", "2025-01-08", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=27")
    ,@("[융합트랙] 2025학년도 스마트인프라SW 융합트랙 이수 신청 안내(학과: 스마트인프라공학과)", "관리자", "2025-01-08", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=26")
    ,@("[융합트랙] 2025학년도 스마트시티&모빌리티 융합트랙 이수 신청 안내(학과: 도시‧교통공학과)", "관리자", "2025-01-08", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=25")
    ,@("[융합트랙] 2025학년도 AI로봇설계 융합트랙 이수 신청 안내(학과: 기계자동차공학부 기계설계공학전공)", "관리자", "2025-01-08", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=24")
    ,@("[연계전공] 2025학년도 지능형스마트팜연계전공 이수 신청 안내(참여학과:스마트팜공학과/주관학과: 소프트웨어학과)", "관리자", "2025-01-08", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=23")
    ,@("[연계전공] 2025학년도 미래자동차SW연계전공 이수 신청 안내(참여학과:미래자동차공학과/주관학과: 스마트정보기술공학과)", "관리자", "2025-01-08", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=22")
    ,@("[연계전공] 2025학년도 AI기반의료정보분석 연계전공 이수 신청 안내(참여학과:의료정보학과/주관학과: 인공지능학부)", "관리자", "2025-01-07", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=21")
    ,@("[홍보] [이스트소프트] 채용연계 교육 과정 (이스트캠프 백엔드, 프론트엔드, iOS 개발자 과정)", "관리자", "2025-01-02", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=20")
    ,@("2024년 SW중심대학사업 SW마일리지 장학금 지급 명단 안내 (총 205명)", "관리자", "2024-12-27", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=19")
    ,@("[추천채용] (주)머신앤비전", "관리자", "2024-12-24", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=18")
    ,@("[홍보] [네이버클라우드] 2024 특별 커리어 성장 세미나 및 정보보안 교육과정 설명회 개최", "관리자", "2024-12-11", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=15")
    ,@("[홍보] LG Aimers 6기 모집 (접수기간 : ~24. 12. 19.(목) 18시까지)", "관리자", "2024-12-11", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=14")
    ,@("[충남테크노파크] 2024년 혁신네트워크 활성화를 위한 청년 구직 희망 조사", "관리자", "2024-12-03", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=10")
    ,@("[소중협] SW중심대학 백서발간을 위한 설문조사(추첨 300명 기프티콘)", "관리자", "2024-12-03", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=9")
    ,@("2024년도 SW인재페스티벌 참가", "관리자", "2024-11-28", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=8")
    ,@("2024년 SW융합캡스톤디자인 경진대회 개최", "관리자", "2024-11-28", "https://swknu.kongju.ac.kr/community/noticedetail.do?seq=7")
)

$startRow = 14
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).NumberFormat = "@"
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
    $ws.Cells.Item($row, 3).ClearFormats()
    $ws.Cells.Item($row, 4).Value = $data[$i][3]
}